$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 7 (ID=6, Greek-text line) -- rows below shift up,
# and the now-unused shared string gets dropped from the workbook on save.
$ws.Rows("7").Delete()

# Update the selection to match the target state.
$ws.Range("F10").Select()
